$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.400803565979004
$ws.Range("B1").Value = 2.230364322662354
$ws.Range("C1").Value = 1.61673104763031
$ws.Range("D1").Value = 1.722981214523315
$ws.Range("E1").Value = 1.586771488189697
